$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Column AW
$ws.Range("AW1").Value = "Please describe how your program serves historicaly underrepresented populations. What actions are you taking to make this an experience for all students to thrive? (Select and briefly describe all that apply.)"
$ws.Range("AW1").Font.Color = 255
$ws.Range("AW2").Value = "Rural"
$ws.Range("AW3").Value = "MOSS instructors are trained in best practices for inclusive education -- honoring diverse perspectives, using student-centered approaches, and acknowledging that students come from different cultural and value frameworks. "

# Column AX
$ws.Range("AX2").Value = "English language learners"
$ws.Range("AX3").Value = "We incorporate visual and hands-on learning as well as oral and written word in our curriculum so ELL students can experience the curriculum through multiple entry points. "

# Column AY
$ws.Range("AY2").Value = "Special education"
$ws.Range("AY3").Value = "We work to accommodate individual student needs, though we rely on the school to provide specifcs regarding IEPs or instructional aids. "

# Column AZ
$ws.Range("AZ2").Value = "Learners with disabilities"
$ws.Range("AZ3").Value = "We work to accommodate individual student needs, though we rely on the school to provide specifcs regarding IEPs or instructional aids. We have ADA accessible bunkhouse units, classrooms, and shower facilities."

# Column BA
$ws.Range("BA2").Value = "Low socio-economic"
$ws.Range("BA3").Value = "We make every effort to make our programs financially accessible. We provide equipment and clothing for students who may not have proper winter clothing."

# Column BB
$ws.Range("BB2").Value = "American Indian/Alaska Native"
$ws.Range("BB3").Value = "See above under rural section "

# Column BC
$ws.Range("BC2").Value = "Asian"
$ws.Range("BC3").Value = "See above under rural section "

# Column BD
$ws.Range("BD2").Value = "Native Hawaiian/Pacific Islander"
$ws.Range("BD3").Value = "See above under rural section "

# Column BE
$ws.Range("BE2").Value = "Black/African American"
$ws.Range("BE3").Value = "See above under rural section "

# Column BF
$ws.Range("BF2").Value = "Hispanic/Latino"
$ws.Range("BF3").Value = "Our enrollment paperwork is available in Spanish as well as English. See above under rural section for other strategies."

# Column BG
$ws.Range("BG2").Value = "Other (list)"
$ws.Range("BG3").Value = "abc"

# Column BH
$ws.Range("BH1").Value = "Regarding the previous question, what efforts are you making to ensure all of your students (including those from historically underrepresented populations) participate in Outdoor School?"
$ws.Range("BH1").Font.Color = 255
$ws.Range("BH3").Value = "All students are able to attend Outdoor School through generous donations and school sponsored fundraisers."

# Column BI
$ws.Range("BI1").Value = "If you requested and received finding for “extenuating circumstances,” please account for how those funds were used.  Be sure to include documentation."
$ws.Range("BI1").Font.Color = 255
$ws.Range("BI3").Value = "n/a"

# Column BJ
$ws.Range("BJ1").Value = "What local (district, community, partner, parent, etc.) resources or funds have you accessed in addition to your state ODS funding?"
$ws.Range("BJ1").Font.Color = 255
$ws.Range("BJ3").Value = "In the past, the teacher representing the class that is attending has budjeted for part of Outdoor School with the rest of the funding coming through fundraisers and individual donations."
